# Updates cryptos list cell values/percentages per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.015.87"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "3.451.57"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'577.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "'148.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").Value = "'0.407"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").Value = "4.041.16"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").Value = "'28.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").Value = "3.454.36"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "'0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").Value = "63.046.95"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "'6.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").Value = "'14.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'9.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.78%  "
$ws.Range("D21").Value = "'385.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.560"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").Value = "'74.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "3.583.67"
$ws.Range("D26").Value = "'0.0000115"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.29%  "
$ws.Range("E27").Value = "  -4.84%  "
$ws.Range("D28").Value = "'7.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "'8.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("D31").Value = "'2.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'23.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'1.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.98%  "
$ws.Range("D35").Value = "'5.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").Value = "'1.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.99%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'7.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "'31.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.05%  "
$ws.Range("D39").Value = "'169.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("D40").Value = "3.486.73"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").Value = "'0.0765"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("D44").Value = "'1.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "'4.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("D46").Value = "'1.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "2.579.51"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").Value = "'22.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.62%  "
